$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3402.8823
$ws.Cells.Item(40, 9).Value = 1677.7778
$ws.Cells.Item(40, 11).Value = 1677.7778
$ws.Cells.Item(40, 13).Value = -1502.7778
$ws.Cells.Item(43, 8).Value = 7332.8335
$ws.Cells.Item(43, 10).Value = 8499.5
$ws.Cells.Item(43, 12).Value = 8499.5
$ws.Cells.Item(43, 14).Value = -8637.5
$ws.Cells.Item(86, 8).Value = 83379070
$ws.Cells.Item(86, 9).Value = 58833590
$ws.Cells.Item(86, 11).Value = 58833590
$ws.Cells.Item(86, 13).Value = -58832467
$ws.Cells.Item(89, 8).Value = 83379070
$ws.Cells.Item(89, 9).Value = 58833590
$ws.Cells.Item(89, 11).Value = 294167950
$ws.Cells.Item(89, 13).Value = -294162334
$ws.Cells.Item(92, 8).Value = 15625776
$ws.Cells.Item(92, 9).Value = 19231584
$ws.Cells.Item(92, 10).Value = 610.3333
$ws.Cells.Item(92, 11).Value = 19231584
$ws.Cells.Item(92, 12).Value = 610.3333
$ws.Cells.Item(92, 13).Value = -19230336
$ws.Cells.Item(92, 14).Value = -3106.3333
$ws.Cells.Item(98, 8).Value = 1042.6471
$ws.Cells.Item(98, 9).Value = 581.73334
$ws.Cells.Item(98, 10).Value = 4499.5
$ws.Cells.Item(98, 11).Value = 581.73334
$ws.Cells.Item(98, 12).Value = 4499.5
$ws.Cells.Item(98, 13).Value = 916.26666
$ws.Cells.Item(98, 14).Value = -7495.5
$ws.Cells.Item(103, 8).Value = 2249.5
$ws.Cells.Item(103, 9).Value = 1666
$ws.Cells.Item(103, 11).Value = 4998
$ws.Cells.Item(103, 13).Value = -4412
$ws.Cells.Item(104, 8).Value = 1213.1666
$ws.Cells.Item(104, 10).Value = 1955
$ws.Cells.Item(104, 12).Value = 5865
$ws.Cells.Item(104, 14).Value = -9359
$ws.Cells.Item(122, 8).Value = 1042.6471
$ws.Cells.Item(122, 9).Value = 581.73334
$ws.Cells.Item(122, 10).Value = 4499.5
$ws.Cells.Item(122, 11).Value = 1745.20002
$ws.Cells.Item(122, 12).Value = 13498.5
$ws.Cells.Item(122, 13).Value = 704.79998
$ws.Cells.Item(122, 14).Value = -18398.5
$ws.Cells.Item(129, 8).Value = 2116.3333
$ws.Cells.Item(129, 9).Value = 1261.875
$ws.Cells.Item(129, 11).Value = 3785.625
$ws.Cells.Item(129, 13).Value = 1214.375
$ws.Cells.Item(131, 8).Value = 2146.6924
$ws.Cells.Item(131, 9).Value = 667.6667
$ws.Cells.Item(131, 10).Value = 5474.5
$ws.Cells.Item(131, 11).Value = 2003.0001
$ws.Cells.Item(131, 12).Value = 16423.5
$ws.Cells.Item(131, 13).Value = 3036.9999
$ws.Cells.Item(131, 14).Value = -26503.5
$ws.Cells.Item(138, 8).Value = 1907.97
$ws.Cells.Item(138, 9).Value = 946.7027
$ws.Cells.Item(138, 10).Value = 2472.524
$ws.Cells.Item(138, 11).Value = 2840.1081
$ws.Cells.Item(138, 12).Value = 7417.572
$ws.Cells.Item(138, 13).Value = 2299.8919
$ws.Cells.Item(138, 14).Value = -17697.572
$ws.Cells.Item(141, 8).Value = 400
$ws.Cells.Item(141, 10).Value = 400
$ws.Cells.Item(141, 12).Value = 1200
$ws.Cells.Item(141, 14).Value = -11560

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22831790
$ws.Cells.Item(32, 9).Value = 30130844
$ws.Cells.Item(32, 11).Value = 30130844
$ws.Cells.Item(32, 13).Value = -30130557
$ws.Cells.Item(45, 8).Value = 2452.8125
$ws.Cells.Item(45, 9).Value = 1284.75
$ws.Cells.Item(45, 11).Value = 1284.75
$ws.Cells.Item(45, 13).Value = -907.75
$ws.Cells.Item(61, 8).Value = 3746.6943
$ws.Cells.Item(61, 9).Value = 3853.682
$ws.Cells.Item(61, 11).Value = 3853.682
$ws.Cells.Item(61, 13).Value = -3641.682
$ws.Cells.Item(102, 8).Value = 2825.375
$ws.Cells.Item(102, 10).Value = 3749.5
$ws.Cells.Item(102, 12).Value = 3749.5
$ws.Cells.Item(102, 14).Value = -6993.5
$ws.Cells.Item(122, 8).Value = 3717.4211
$ws.Cells.Item(122, 9).Value = 2438
$ws.Cells.Item(122, 10).Value = 4868.9
$ws.Cells.Item(122, 11).Value = 7314
$ws.Cells.Item(122, 12).Value = 14606.7
$ws.Cells.Item(122, 13).Value = -4864
$ws.Cells.Item(122, 14).Value = -19506.7
$ws.Cells.Item(132, 8).Value = 225432.36
$ws.Cells.Item(132, 9).Value = 296352.9
$ws.Cells.Item(132, 11).Value = 889058.7000000001
$ws.Cells.Item(132, 13).Value = -886528.7000000001
$ws.Cells.Item(136, 8).Value = 3746.6943
$ws.Cells.Item(136, 9).Value = 3853.682
$ws.Cells.Item(136, 11).Value = 11561.046
$ws.Cells.Item(136, 13).Value = -9011.045999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2683.3
$ws.Cells.Item(107, 9).Value = 2322.5
$ws.Cells.Item(107, 10).Value = 3224.5
$ws.Cells.Item(107, 11).Value = 2322.5
$ws.Cells.Item(107, 12).Value = 3224.5
$ws.Cells.Item(107, 13).Value = -402.5
$ws.Cells.Item(107, 14).Value = -7064.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5136.9395
$ws.Cells.Item(31, 9).Value = 2208.6765
$ws.Cells.Item(31, 11).Value = 2208.6765
$ws.Cells.Item(31, 13).Value = -1913.6765
$ws.Cells.Item(34, 8).Value = 5136.9395
$ws.Cells.Item(34, 9).Value = 2208.6765
$ws.Cells.Item(34, 11).Value = 2208.6765
$ws.Cells.Item(34, 13).Value = -2006.6765
$ws.Cells.Item(99, 8).Value = 3276.5557
$ws.Cells.Item(99, 9).Value = 3081.8333
$ws.Cells.Item(99, 10).Value = 3666
$ws.Cells.Item(99, 11).Value = 3081.8333
$ws.Cells.Item(99, 12).Value = 3666
$ws.Cells.Item(99, 13).Value = -1583.8333
$ws.Cells.Item(99, 14).Value = -6662
$ws.Cells.Item(126, 8).Value = 3276.5557
$ws.Cells.Item(126, 9).Value = 3081.8333
$ws.Cells.Item(126, 10).Value = 3666
$ws.Cells.Item(126, 11).Value = 9245.499899999999
$ws.Cells.Item(126, 12).Value = 10998
$ws.Cells.Item(126, 13).Value = -6775.499899999999
$ws.Cells.Item(126, 14).Value = -15938
$ws.Cells.Item(132, 8).Value = 4486.9644
$ws.Cells.Item(132, 9).Value = 4432.143
$ws.Cells.Item(132, 10).Value = 4651.4287
$ws.Cells.Item(132, 11).Value = 13296.429
$ws.Cells.Item(132, 12).Value = 13954.2861
$ws.Cells.Item(132, 13).Value = -10766.429
$ws.Cells.Item(132, 14).Value = -19014.2861
$ws.Cells.Item(134, 8).Value = 2710
$ws.Cells.Item(134, 9).Value = 2761.5
$ws.Cells.Item(134, 11).Value = 8284.5
$ws.Cells.Item(134, 13).Value = -5749.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 4845.7
$ws.Cells.Item(18, 9).Value = 4076.3333
$ws.Cells.Item(18, 11).Value = 12228.9999
$ws.Cells.Item(18, 13).Value = -12059.9999
$ws.Cells.Item(34, 8).Value = 648.875
$ws.Cells.Item(34, 9).Value = 87.166664
$ws.Cells.Item(34, 10).Value = 2334
$ws.Cells.Item(34, 11).Value = 261.499992
$ws.Cells.Item(34, 12).Value = 7002
$ws.Cells.Item(34, 13).Value = -177.499992
$ws.Cells.Item(34, 14).Value = -7170
$ws.Cells.Item(39, 8).Value = 4903.1
$ws.Cells.Item(39, 10).Value = 4903.1
$ws.Cells.Item(39, 12).Value = 14709.3
$ws.Cells.Item(39, 14).Value = -15297.3
$ws.Cells.Item(55, 8).Value = 3113.0588
$ws.Cells.Item(55, 9).Value = 985.875
$ws.Cells.Item(55, 10).Value = 5003.8887
$ws.Cells.Item(55, 11).Value = 2957.625
$ws.Cells.Item(55, 12).Value = 15011.6661
$ws.Cells.Item(55, 13).Value = -2780.625
$ws.Cells.Item(55, 14).Value = -15365.6661
$ws.Cells.Item(132, 8).Value = 1098.5714
$ws.Cells.Item(132, 10).Value = 990
$ws.Cells.Item(132, 12).Value = 8910
$ws.Cells.Item(132, 14).Value = -13970

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 24002.5
$ws.Cells.Item(18, 9).Value = 24002.5
$ws.Cells.Item(18, 11).Value = 24002.5
$ws.Cells.Item(18, 13).Value = -23709.5
$ws.Cells.Item(70, 8).Value = 22194
$ws.Cells.Item(70, 9).Value = 82797
$ws.Cells.Item(70, 10).Value = 6033.2
$ws.Cells.Item(70, 11).Value = 82797
$ws.Cells.Item(70, 12).Value = 6033.2
$ws.Cells.Item(70, 13).Value = -82527
$ws.Cells.Item(70, 14).Value = -6573.2
$ws.Cells.Item(73, 8).Value = 22194
$ws.Cells.Item(73, 9).Value = 82797
$ws.Cells.Item(73, 10).Value = 6033.2
$ws.Cells.Item(73, 11).Value = 82797
$ws.Cells.Item(73, 12).Value = 6033.2
$ws.Cells.Item(73, 13).Value = -81861
$ws.Cells.Item(73, 14).Value = -7905.2
$ws.Cells.Item(113, 8).Value = 16421.715
$ws.Cells.Item(113, 9).Value = 14646.375
$ws.Cells.Item(113, 10).Value = 18788.834
$ws.Cells.Item(113, 11).Value = 14646.375
$ws.Cells.Item(113, 12).Value = 18788.834
$ws.Cells.Item(113, 13).Value = -12476.375
$ws.Cells.Item(113, 14).Value = -23128.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2830.75
$ws.Cells.Item(22, 9).Value = 1189.2
$ws.Cells.Item(22, 10).Value = 5566.6665
$ws.Cells.Item(22, 11).Value = 1189.2
$ws.Cells.Item(22, 12).Value = 5566.6665
$ws.Cells.Item(22, 13).Value = -894.2
$ws.Cells.Item(22, 14).Value = -6156.6665
$ws.Cells.Item(27, 8).Value = 2830.75
$ws.Cells.Item(27, 9).Value = 1189.2
$ws.Cells.Item(27, 10).Value = 5566.6665
$ws.Cells.Item(27, 11).Value = 1189.2
$ws.Cells.Item(27, 12).Value = 5566.6665
$ws.Cells.Item(27, 13).Value = -1082.2
$ws.Cells.Item(27, 14).Value = -5780.6665
$ws.Cells.Item(61, 8).Value = 1369.8
$ws.Cells.Item(61, 9).Value = 966.44446
$ws.Cells.Item(61, 11).Value = 966.44446
$ws.Cells.Item(61, 13).Value = -764.44446
$ws.Cells.Item(113, 8).Value = 1369.8
$ws.Cells.Item(113, 9).Value = 966.44446
$ws.Cells.Item(113, 11).Value = 966.44446
$ws.Cells.Item(113, 13).Value = 1203.55554

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 11005
$ws.Cells.Item(19, 9).Value = 11005
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 11005
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -10831
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(54, 8).Value = 30070
$ws.Cells.Item(54, 9).Value = 30070
$ws.Cells.Item(54, 11).Value = 30070
$ws.Cells.Item(54, 13).Value = -29550
$ws.Cells.Item(62, 8).Value = 5715.5
$ws.Cells.Item(62, 10).Value = 6208.6
$ws.Cells.Item(62, 12).Value = 6208.6
$ws.Cells.Item(62, 14).Value = -7456.6
$ws.Cells.Item(65, 8).Value = 5715.5
$ws.Cells.Item(65, 10).Value = 6208.6
$ws.Cells.Item(65, 12).Value = 31043
$ws.Cells.Item(65, 14).Value = -37283
$ws.Cells.Item(81, 8).Value = 70997.47
$ws.Cells.Item(81, 10).Value = 5628.2856
$ws.Cells.Item(81, 12).Value = 11256.5712
$ws.Cells.Item(81, 14).Value = -13378.5712
$ws.Cells.Item(84, 8).Value = 70997.47
$ws.Cells.Item(84, 10).Value = 5628.2856
$ws.Cells.Item(84, 12).Value = 56282.856
$ws.Cells.Item(84, 14).Value = -66890.856
$ws.Cells.Item(135, 8).Value = 99999
$ws.Cells.Item(135, 10).Value = 99999
$ws.Cells.Item(135, 12).Value = 99999
$ws.Cells.Item(135, 14).Value = -110139
